# Reposition the decorative "Group 7" shape on slide 1.
#
# The underlying OOXML moved this group's <a:xfrm><a:off> from (1, 0) EMU to
# (139700, 153279) EMU while leaving its extent (<a:ext>) and child offset/
# extent (<a:chOff>/<a:chExt>) untouched. PowerPoint's Shape.Left / Shape.Top
# properties are expressed in points (1 pt = 12700 EMU), so convert.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

$targetShape = $null
for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $sh = $s.Shapes.Item($i)
    if ($sh.Name -eq "Group 7") {
        $targetShape = $sh
        break
    }
}

if ($targetShape -eq $null) {
    $targetShape = $s.Shapes.Item(2)
}

$emuPerPoint = 12700
$targetShape.Left = 139700 / $emuPerPoint
$targetShape.Top = 153279 / $emuPerPoint
